$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'37.466.16"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.72%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.077.35"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.09%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.07%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'232.63"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.13%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'0.631"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +1.31%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.06%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'57.58"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -1.38%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.389"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -1.02%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.0779"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -0.53%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  +2.29%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'14.94"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +1.44%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'2.386.87"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.21%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'20.95"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +0.57%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.769"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -0.34%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'5.33"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.11%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'2.073.53"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -0.12%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'37.433.02"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.59%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("B19").Value = "'Litecoin"
$ws.Range("B19").Style = "Normal"
$ws.Range("C19").Value = "'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("C19").Style = "Normal"
$ws.Range("D19").Value = "'70.55"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -0.70%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("B20").Value = "'Uniswap"
$ws.Range("B20").Style = "Normal"
$ws.Range("C20").Value = "'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("C20").Style = "Normal"
$ws.Range("D20").Value = "'6.02"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -2.59%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.0₃0830"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -0.31%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'228.53"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.23%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.02%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'2.37"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -0.64%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'2.33"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -2.33%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'9.64"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +6.87%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'170.08"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -0.14%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  -3.49%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'19.53"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +0.42%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'1.37"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -1.16%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  +1.71%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'4.63"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -0.98%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.0632"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +0.32%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'4.63"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -0.37%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'2.47"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -0.35%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'1.81"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -0.16%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'3.30"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -2.65%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  +0.04%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'5.28"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -0.74%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.0231"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +7.25%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'100.03"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -0.46%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = "'HuobiToken"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'2.92"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +1.25%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("B43").Value = "'TrustWalletToken"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = "'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'1.20"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +4.13%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("B44").Value = "'Cronos"
$ws.Range("B44").Style = "Normal"
$ws.Range("C44").Value = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").Value = "'0.0952"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -1.79%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("B45").Value = "'Maker"
$ws.Range("B45").Style = "Normal"
$ws.Range("C45").Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").Value = "'1.463.77"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +1.27%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("B46").Value = "'InjectiveProtocol"
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = "'16.68"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +2.09%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'1.04"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -1.25%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("B48").Value = "'FraxShare"
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value = "'7.24"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -2.20%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("B49").Value = "'FTXToken"
$ws.Range("B49").Style = "Normal"
$ws.Range("C49").Value = "'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("C49").Style = "Normal"
$ws.Range("D49").Value = "'3.93"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -6.31%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'2.94"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -2.09%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'2.270.96"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +0.17%  "
$ws.Range("E51").Style = "Normal"
